$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Append new BOM rows (58-61) ----------------------------------------
# Values are written in this specific order so new shared-string entries
# land at the same indices the original author produced.
$ws.Range("C58").Value = "CONN AUDIO JACK 3.5MM 4COND SMD"
$ws.Range("D58").Value = "SJ-43516-SMT-TR"
$ws.Range("E58").Value = "CP-43516SJCT-ND"

$ws.Range("C59").Value = "CONN JACK 3.5MM R/A 4POS MID SMD"

$ws.Range("B58").Value = "0"

$ws.Range("D59").Value = "SJ-43617-SMT-TR"
$ws.Range("E59").Value = "CP-43617SJCT-ND"

$ws.Range("C60").Value = "IC REG LDO 3.3V 0.3A SOT23-5"
$ws.Range("D60").Value = "MCP1802T-3302I/OT"
$ws.Range("E60").Value = "MCP1802T-3302I/OTCT-ND"

$ws.Range("C61").Value = "IC POT DGTL 256-TAP 10UMAX"
$ws.Range("D61").Value = "MAX5388NAUB+"
$ws.Range("E61").Value = "MAX5388NAUB+-ND"

$ws.Range("B59").Value = "1"
$ws.Range("B60").Value = "1"
$ws.Range("B61").Value = "1"

# --- Update view / selection state --------------------------------------
$ws.Activate()
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 2
$ws.Range("E64").Select()
